# Applies the cryptos.xlsx data refresh for Tue May  2 06:51:39 UTC 2023 (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.015.52'
$ws.Range("E2").Value = '  -2.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.829.50'
$ws.Range("E3").Value = '  -1.14%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.40'
$ws.Range("E5").Value = '  -2.88%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4651'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3861'
$ws.Range("E8").Value = '  -1.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07856'
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9577'
$ws.Range("E10").Value = '  -2.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.86'
$ws.Range("E11").Value = '  -1.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.844.57'
$ws.Range("E12").Value = '  -10.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.675'
$ws.Range("E13").Value = '  -3.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.903'
$ws.Range("E14").Value = '  -1.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06860'
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("E16").Value = '  -0.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009909'
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("E19").Value = '  -3.35%  '
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '28.036.77'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.310'
$ws.Range("E22").Value = '  -1.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.97'
$ws.Range("E23").Value = '  -2.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.096'
$ws.Range("E24").Value = '  -1.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.124.57'
$ws.Range("E25").Value = '  -6.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.53'
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.08'
$ws.Range("E27").Value = '  -1.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.659'
$ws.Range("E28").Value = '  -7.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.955'
$ws.Range("E29").Value = '  -3.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.36'
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09227'
$ws.Range("E31").Value = '  -2.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9306'
$ws.Range("E32").Value = '  -5.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.261'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.317'
$ws.Range("E34").Value = '  -2.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.294'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05834'
$ws.Range("E36").Value = '  -5.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02114'
$ws.Range("E37").Value = '  -3.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.133'
$ws.Range("E38").Value = '  -2.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.794'
$ws.Range("E39").Value = '  +2.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5582'
$ws.Range("E40").Value = '  -2.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.862'
$ws.Range("E41").Value = '  -2.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1757'
$ws.Range("E42").Value = '  -2.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.07185'
$ws.Range("E43").Value = '  +0.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.56'
$ws.Range("E44").Value = '  -2.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5255'
$ws.Range("E45").Value = '  -2.59%  '
$ws.Range("E46").Value = '  -10.94%  '
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.118'
$ws.Range("E47").Value = '  -10.47%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.826'
$ws.Range("E48").Value = '  -4.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '112.45'
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.327'
$ws.Range("E51").Value = '  +0.52%  '
